$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptn"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3206776666666667
$ws.Range("H2").Value = 0.962033
$ws.Range("I2").Value = 0.03001977461414601
$ws.Range("J2").Value = 0.03001977461414601
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03995766666666666
$ws.Range("N2").Value = 0.119873
$ws.Range("O2").Value = 0.005314930928687666
$ws.Range("P2").Value = 0.005314930928687667
$ws.Range("Q2").Value = 0.01281353131211111
$ws.Range("R2").Value = 0.115321781809
$ws.Range("S2").Value = 0.0001595530285689575
$ws.Range("T2").Value = 0.0001595530285689575

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptn"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3206776666666667
$ws.Range("H3").Value = 0.962033
$ws.Range("I3").Value = 0.03001977461414601
$ws.Range("J3").Value = 0.03001977461414601
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03069133333333333
$ws.Range("N3").Value = 0.092074
$ws.Range("O3").Value = 0.004082378436578614
$ws.Range("P3").Value = 0.004082378436578615
$ws.Range("Q3").Value = 0.009842025160222223
$ws.Range("R3").Value = 0.088578226442
$ws.Range("S3").Value = 0.0001225520805557398
$ws.Range("T3").Value = 0.0001225520805557398

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ptn"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3206776666666667
$ws.Range("H4").Value = 0.962033
$ws.Range("I4").Value = 0.03001977461414601
$ws.Range("J4").Value = 0.03001977461414601
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.447354000000001
$ws.Range("N4").Value = 22.342062
$ws.Range("O4").Value = 0.9906026906347337
$ws.Range("P4").Value = 0.9906026906347338
$ws.Range("Q4").Value = 2.388200103560667
$ws.Range("R4").Value = 21.493800932046
$ws.Range("S4").Value = 0.02973766950502131
$ws.Range("T4").Value = 0.02973766950502132

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ptn"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.752274333333333
$ws.Range("H5").Value = 17.256823
$ws.Range("I5").Value = 0.5384908178993973
$ws.Range("J5").Value = 0.5384908178993975
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03995766666666666
$ws.Range("N5").Value = 0.119873
$ws.Range("O5").Value = 0.005314930928687666
$ws.Range("P5").Value = 0.005314930928687667
$ws.Range("Q5").Value = 0.2298474603865555
$ws.Range("R5").Value = 2.068627143479
$ws.Range("S5").Value = 0.002862041502867825
$ws.Range("T5").Value = 0.002862041502867826

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ptn"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.752274333333333
$ws.Range("H6").Value = 17.256823
$ws.Range("I6").Value = 0.5384908178993973
$ws.Range("J6").Value = 0.5384908178993975
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03069133333333333
$ws.Range("N6").Value = 0.092074
$ws.Range("O6").Value = 0.004082378436578614
$ws.Range("P6").Value = 0.004082378436578615
$ws.Range("Q6").Value = 0.1765449689891111
$ws.Range("R6").Value = 1.588904720902
$ws.Range("S6").Value = 0.002198323303288081
$ws.Range("T6").Value = 0.002198323303288082

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ptn"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.752274333333333
$ws.Range("H7").Value = 17.256823
$ws.Range("I7").Value = 0.5384908178993973
$ws.Range("J7").Value = 0.5384908178993975
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.447354000000001
$ws.Range("N7").Value = 22.342062
$ws.Range("O7").Value = 0.9906026906347337
$ws.Range("P7").Value = 0.9906026906347338
$ws.Range("Q7").Value = 42.83922326544734
$ws.Range("R7").Value = 385.553009389026
$ws.Range("S7").Value = 0.5334304530932414
$ws.Range("T7").Value = 0.5334304530932417

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ptn"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.609262333333334
$ws.Range("H8").Value = 13.827787
$ws.Range("I8").Value = 0.4314894074864565
$ws.Range("J8").Value = 0.4314894074864565
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03995766666666666
$ws.Range("N8").Value = 0.119873
$ws.Range("O8").Value = 0.005314930928687666
$ws.Range("P8").Value = 0.005314930928687667
$ws.Range("Q8").Value = 0.1841753678945555
$ws.Range("R8").Value = 1.657578311051
$ws.Range("S8").Value = 0.002293336397250883
$ws.Range("T8").Value = 0.002293336397250883

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ptn"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.609262333333334
$ws.Range("H9").Value = 13.827787
$ws.Range("I9").Value = 0.4314894074864565
$ws.Range("J9").Value = 0.4314894074864565
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.03069133333333333
$ws.Range("N9").Value = 0.092074
$ws.Range("O9").Value = 0.004082378436578614
$ws.Range("P9").Value = 0.004082378436578615
$ws.Range("Q9").Value = 0.1414644066931111
$ws.Range("R9").Value = 1.273179660238
$ws.Range("S9").Value = 0.001761503052734793
$ws.Range("T9").Value = 0.001761503052734793

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ptn"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.609262333333334
$ws.Range("H10").Value = 13.827787
$ws.Range("I10").Value = 0.4314894074864565
$ws.Range("J10").Value = 0.4314894074864565
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.447354000000001
$ws.Range("N10").Value = 22.342062
$ws.Range("O10").Value = 0.9906026906347337
$ws.Range("P10").Value = 0.9906026906347338
$ws.Range("Q10").Value = 34.32680827519934
$ws.Range("R10").Value = 308.941274476794
$ws.Range("S10").Value = 0.4274345680364708
$ws.Range("T10").Value = 0.4274345680364709
